$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.066.36'
$ws.Range("E2").Value = '  +3.47%  '
$ws.Range("D3").Value = '3.194.99'
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '538.05'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").Value = '145.07'
$ws.Range("E6").Value = '  +4.24%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +3.52%  '
$ws.Range("D9").Value = '7.35'
$ws.Range("D10").Value = '0.114'
$ws.Range("E10").Value = '  +4.29%  '
$ws.Range("D11").Value = '0.431'
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("D12").Value = '3.745.85'
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").Value = '0.0000175'
$ws.Range("E14").Value = '  +3.04%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '26.06'
$ws.Range("E15").Value = '  +1.21%  '
$ws.Range("D16").Value = '60.101.21'
$ws.Range("E16").Value = '  +3.42%  '
$ws.Range("D17").Value = '3.197.62'
$ws.Range("E17").Value = '  +1.93%  '
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("D19").Value = '13.09'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = '8.36'
$ws.Range("E20").Value = '  +1.50%  '
$ws.Range("D21").Value = '384.87'
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '0.530'
$ws.Range("E23").Value = '  +2.56%  '
$ws.Range("D24").Value = '70.38'
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("E25").Value = '  +2.34%  '
$ws.Range("D26").Value = '8.86'
$ws.Range("E26").Value = '  +11.40%  '
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("D28").Value = '0.0₃0905'
$ws.Range("E28").Value = '  +2.18%  '
$ws.Range("D29").Value = '1.91'
$ws.Range("E29").Value = '  +1.04%  '
$ws.Range("D30").Value = '22.42'
$ws.Range("E30").Value = '  +3.04%  '
$ws.Range("D31").Value = '6.18'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").Value = '5.37'
$ws.Range("E32").Value = '  +3.73%  '
$ws.Range("E33").Value = '  +4.01%  '
$ws.Range("D34").Value = '6.58'
$ws.Range("E34").Value = '  +4.57%  '
$ws.Range("D35").Value = '156.07'
$ws.Range("E35").Value = '  -3.37%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '2.775.50'
$ws.Range("E37").Value = '  +5.43%  '
$ws.Range("D38").Value = '25.78'
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("D39").Value = '0.0715'
$ws.Range("E39").Value = '  +5.92%  '
$ws.Range("D40").Value = '1.69'
$ws.Range("E40").Value = '  +0.97%  '
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").Value = '39.76'
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("D43").Value = '0.730'
$ws.Range("E43").Value = '  +4.26%  '
$ws.Range("D44").Value = '0.0286'
$ws.Range("E44").Value = '  +5.07%  '
$ws.Range("D45").Value = '3.234.37'
$ws.Range("E45").Value = '  +1.75%  '
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +2.60%  '
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("D48").Value = '6.18'
$ws.Range("D49").Value = '0.803'
$ws.Range("E49").Value = '  +7.17%  '
$ws.Range("D50").Value = '20.57'
$ws.Range("E51").Value = '  +0.00%  '
